$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the worksheet dimension-relevant data: rows 2-9 are updated with recomputed
# NATMI values, and new rows 10-13 are appended for the additional "sCs" sending cluster.
$arr = New-Object "object[,]" 12,20

# Row 2 (array row 0)
$arr[0,0] = "ECs"
$arr[0,1] = "Bmp6"
$arr[0,2] = "Acvr1"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 19.72450833333333
$arr[0,7] = 59.173525
$arr[0,8] = 0.5834853563809828
$arr[0,9] = 0.5834853563809829
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.436778333333334
$arr[0,13] = 16.310335
$arr[0,14] = 0.121853993972124
$arr[0,15] = 0.121853993972124
$arr[0,16] = 107.2377795423194
$arr[0,17] = 965.1400158808751
$arr[0,18] = 0.07110002109927091
$arr[0,19] = 0.07110002109927094

# Row 3 (array row 1)
$arr[1,0] = "ECs"
$arr[1,1] = "Bmp6"
$arr[1,2] = "Acvr1"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 19.72450833333333
$arr[1,7] = 59.173525
$arr[1,8] = 0.5834853563809828
$arr[1,9] = 0.5834853563809829
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 20.81735866666667
$arr[1,13] = 62.452076
$arr[1,14] = 0.4665774732677552
$arr[1,15] = 0.4665774732677551
$arr[1,16] = 410.6121644986555
$arr[1,17] = 3695.5094804879
$arr[1,18] = 0.2722411232689746
$arr[1,19] = 0.2722411232689746

# Row 4 (array row 2)
$arr[2,0] = "ECs"
$arr[2,1] = "Bmp6"
$arr[2,2] = "Acvr1"
$arr[2,3] = "M2"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 19.72450833333333
$arr[2,7] = 59.173525
$arr[2,8] = 0.5834853563809828
$arr[2,9] = 0.5834853563809829
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 5.131792
$arr[2,13] = 15.395376
$arr[2,14] = 0.1150183643869107
$arr[2,15] = 0.1150183643869107
$arr[2,16] = 101.2220740689333
$arr[2,17] = 910.9986666204001
$arr[2,18] = 0.06711153133465432
$arr[2,19] = 0.06711153133465433

# Row 5 (array row 3)
$arr[3,0] = "ECs"
$arr[3,1] = "Bmp6"
$arr[3,2] = "Acvr1"
$arr[3,3] = "sCs"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 19.72450833333333
$arr[3,7] = 59.173525
$arr[3,8] = 0.5834853563809828
$arr[3,9] = 0.5834853563809829
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 13.23122433333333
$arr[3,13] = 39.693673
$arr[3,14] = 0.2965501683732101
$arr[3,15] = 0.2965501683732101
$arr[3,16] = 260.9793946230361
$arr[3,17] = 2348.814551607325
$arr[3,18] = 0.1730326806780829
$arr[3,19] = 0.1730326806780829

# Row 6 (array row 4)
$arr[4,0] = "FAPs"
$arr[4,1] = "Bmp6"
$arr[4,2] = "Acvr1"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.7684289999999999
$arr[4,7] = 2.305287
$arr[4,8] = 0.02273146997336134
$arr[4,9] = 0.02273146997336134
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 5.436778333333334
$arr[4,13] = 16.310335
$arr[4,14] = 0.121853993972124
$arr[4,15] = 0.121853993972124
$arr[4,16] = 4.177778137904999
$arr[4,17] = 37.600003241145
$arr[4,18] = 0.002769920405111491
$arr[4,19] = 0.002769920405111492

# Row 7 (array row 5)
$arr[5,0] = "FAPs"
$arr[5,1] = "Bmp6"
$arr[5,2] = "Acvr1"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 0.7684289999999999
$arr[5,7] = 2.305287
$arr[5,8] = 0.02273146997336134
$arr[5,9] = 0.02273146997336134
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 20.81735866666667
$arr[5,13] = 62.452076
$arr[5,14] = 0.4665774732677552
$arr[5,15] = 0.4665774732677551
$arr[5,16] = 15.996662102868
$arr[5,17] = 143.969958925812
$arr[5,18] = 0.01060599182383278
$arr[5,19] = 0.01060599182383278

# Row 8 (array row 6)
$arr[6,0] = "FAPs"
$arr[6,1] = "Bmp6"
$arr[6,2] = "Acvr1"
$arr[6,3] = "M2"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 0.7684289999999999
$arr[6,7] = 2.305287
$arr[6,8] = 0.02273146997336134
$arr[6,9] = 0.02273146997336134
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 5.131792
$arr[6,13] = 15.395376
$arr[6,14] = 0.1150183643869107
$arr[6,15] = 0.1150183643869107
$arr[6,16] = 3.943417794768
$arr[6,17] = 35.490760152912
$arr[6,18] = 0.002614536496446193
$arr[6,19] = 0.002614536496446194

# Row 9 (array row 7)
$arr[7,0] = "FAPs"
$arr[7,1] = "Bmp6"
$arr[7,2] = "Acvr1"
$arr[7,3] = "sCs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 0.7684289999999999
$arr[7,7] = 2.305287
$arr[7,8] = 0.02273146997336134
$arr[7,9] = 0.02273146997336134
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 13.23122433333333
$arr[7,13] = 39.693673
$arr[7,14] = 0.2965501683732101
$arr[7,15] = 0.2965501683732101
$arr[7,16] = 10.167256483239
$arr[7,17] = 91.505308349151
$arr[7,18] = 0.006741021247970875
$arr[7,19] = 0.006741021247970875

# Row 10 (array row 8)
$arr[8,0] = "sCs"
$arr[8,1] = "Bmp6"
$arr[8,2] = "Acvr1"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 13.31169566666667
$arr[8,7] = 39.935087
$arr[8,8] = 0.3937831736456558
$arr[8,9] = 0.3937831736456558
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 5.436778333333334
$arr[8,13] = 16.310335
$arr[8,14] = 0.121853993972124
$arr[8,15] = 0.121853993972124
$arr[8,16] = 72.37273858046056
$arr[8,17] = 651.3546472241451
$arr[8,18] = 0.04798405246774161
$arr[8,19] = 0.04798405246774162

# Row 11 (array row 9)
$arr[9,0] = "sCs"
$arr[9,1] = "Bmp6"
$arr[9,2] = "Acvr1"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 13.31169566666667
$arr[9,7] = 39.935087
$arr[9,8] = 0.3937831736456558
$arr[9,9] = 0.3937831736456558
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 20.81735866666667
$arr[9,13] = 62.452076
$arr[9,14] = 0.4665774732677552
$arr[9,15] = 0.4665774732677551
$arr[9,16] = 277.1143431545125
$arr[9,17] = 2494.029088390612
$arr[9,18] = 0.1837303581749478
$arr[9,19] = 0.1837303581749478

# Row 12 (array row 10)
$arr[10,0] = "sCs"
$arr[10,1] = "Bmp6"
$arr[10,2] = "Acvr1"
$arr[10,3] = "M2"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 13.31169566666667
$arr[10,7] = 39.935087
$arr[10,8] = 0.3937831736456558
$arr[10,9] = 0.3937831736456558
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 5.131792
$arr[10,13] = 15.395376
$arr[10,14] = 0.1150183643869107
$arr[10,15] = 0.1150183643869107
$arr[10,16] = 68.31285332863467
$arr[10,17] = 614.8156799577121
$arr[10,18] = 0.04529229655581016
$arr[10,19] = 0.04529229655581016

# Row 13 (array row 11)
$arr[11,0] = "sCs"
$arr[11,1] = "Bmp6"
$arr[11,2] = "Acvr1"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 13.31169566666667
$arr[11,7] = 39.935087
$arr[11,8] = 0.3937831736456558
$arr[11,9] = 0.3937831736456558
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 13.23122433333333
$arr[11,13] = 39.693673
$arr[11,14] = 0.2965501683732101
$arr[11,15] = 0.2965501683732101
$arr[11,16] = 176.1300316227279
$arr[11,17] = 1585.170284604551
$arr[11,18] = 0.1167764664471563
$arr[11,19] = 0.1167764664471563

$ws.Range("A2:T13").Value = $arr

